$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 5 new rows at their final positions, using ORIGINAL (pre-insert)
# row numbers and working from the bottom up so earlier insertion points are
# not shifted by later ones.
$ws.Rows(6).Insert()   # new row before old row6 (+5521997432262 ...)
$ws.Rows(5).Insert()   # new row before old row5 (+5511952381413 ...)
$ws.Rows(4).Insert()   # new row before old row4 (+5511967085107 ...)
$ws.Rows(2).Insert()   # two new rows before old row2 (+555491557534 ...)
$ws.Rows(2).Insert()

# Only the two brand-new rows inserted directly under the header (rows 2-3)
# inherit the header's style (bold/red) from Insert(). Every other row keeps
# the normal data-row style it already had. Fix rows 2-3 by copying the
# format from a row that is guaranteed to still have the normal data style.
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A2:C3").PasteSpecial(-4122) | Out-Null

# Force text formatting on the data columns so numeric-looking strings
# (phone numbers, DDD codes, date strings) are stored as plain text, not
# auto-converted to numbers/dates.
$ws.Range("A2:C17").NumberFormat = "@"

$data = @(
    @("+555199199744", "51", "2024-10-10"),
    @("+5522981222545", "22", "2024-10-05"),
    @("+555491557534", "54", "2024-10-04"),
    @("+556181971614", "61", "2024-10-01"),
    @("+5513988453610", "13", "2024-09-25"),
    @("+5511967085107", "11", "2024-09-20"),
    @("+5511977696904", "11", "2024-09-20"),
    @("+5511952381413", "11", "2024-09-18"),
    @("+5521965489343", "21", "2024-09-11"),
    @("+5521997432262", "21", "2024-09-10"),
    @("+5521985109311", "21", "2024-09-09"),
    @("+553291004823", "32", "2024-08-26"),
    @("+34603138909", $null, "2024-08-19"),
    @("+5511967859426", "11", "2024-07-28"),
    @("+5521965197022", "21", "2024-07-21"),
    @("+556298529715", "62", "2024-07-09")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
